$wb = $excel.ActiveWorkbook

# --- Productdata sheet ---
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("C6").Value = 11

# --- Capacity sheet ---
$wsCapacity = $wb.Worksheets.Item("Capacity")
$wsCapacity.Range("B2").Value = 35
$wsCapacity.Range("B3").Value = 175
$wsCapacity.Range("B4").Value = 140
$wsCapacity.Range("B6").Value = 175

# --- ProcessingTime sheet ---
$wsProcessingTime = $wb.Worksheets.Item("ProcessingTime")
$wsProcessingTime.Range("B2").Value = 1
$wsProcessingTime.Range("C3").Value = 5
$wsProcessingTime.Range("D4").Value = 4
$wsProcessingTime.Range("F6").Value = 5
